$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 141, shifting existing rows 141:161 down to 142:162
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with the new record's data
$ws.Cells.Item(141, 1).Value = 5
$ws.Cells.Item(141, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(141, 3).Value = "Maule"
$ws.Cells.Item(141, 4).Value = 44984
$ws.Cells.Item(141, 5).Value = 7
$ws.Cells.Item(141, 6).Value = "Fruta"
$ws.Cells.Item(141, 7).Value = 100108
$ws.Cells.Item(141, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(141, 9).Value = 100108002
$ws.Cells.Item(141, 10).Value = "Mango"
$ws.Cells.Item(141, 11).Value = "Sin especificar"
$ws.Cells.Item(141, 12).Value = "Primera"
$ws.Cells.Item(141, 13).Value = 248
$ws.Cells.Item(141, 14).Value = 8000
$ws.Cells.Item(141, 15).Value = 8000
$ws.Cells.Item(141, 16).Value = 8000
$ws.Cells.Item(141, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(141, 18).Value = "Perú"
$ws.Cells.Item(141, 19).Value = 2000
$ws.Cells.Item(141, 20).Value = 4
